# Update Name of Algo
# Update a handful of numeric values in column A (and D13) of Sheet1
# to reflect the new RandomForest imputation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value  = -21.11050000000001
$ws.Range("A10").Value = -20.58409999999997
$ws.Range("A12").Value = -22.49650000000004
$ws.Range("D13").Value = -7.668700000000005
$ws.Range("A18").Value = -22.35810000000004
$ws.Range("A25").Value = -22.28610000000004
